# Update "Ciudades" sheet with new provincias Spain COVID numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Updated timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 12:35"

# Updated per-row statistics (Casos totales, Casos activos, Recuperados, Muertes).
$updates = @{
    4  = @(65269, 39984, 16525, 8760)
    5  = @(55280, 25326, 24262, 5692)
    6  = @(18105, 7544, 8642, 1919)
    7  = @(16453, 6205, 7413, 2835)
    9  = @(12317, 9481, 1504, 1332)
    11 = @(9301, 7996, 708, 597)
    14 = @(5372, 3312, 1231, 829)
    16 = @(5094, 3297, 1303, 494)
    23 = @(4009, 2804, 857, 348)
    33 = @(2356, 1038, 1011, 307)
    34 = @(2271, 1475, 645, 151)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
